# Workbook "ValueSet-KLConditionCodesNursing" — bump published term to 1.1.0
# (matches commit message "Added 1.1.0 of term").
#
# The "Metadata" sheet lists Property/Value pairs; row 3 holds the Version
# and row 8 holds the Date the value set was published.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Version: 1.0.0 -> 1.1.0
$ws.Range("B3").Value = "1.1.0"

# Date: 2023-06-07T11:52:14+02:00 -> 2023-07-10T23:08:03+02:00
$ws.Range("B8").Value = "2023-07-10T23:08:03+02:00"
